$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fidelity_qubits")

# Ensure the range keeps its text (string) formatting so that values
# round-trip as exact-precision text, matching the original data which
# was stored as inline strings rather than numeric cells.
$ws.Range("A1:E10").NumberFormat = "@"

$ws.Range("A1").Value = "0.5559839614875356"
$ws.Range("B1").Value = "0.9276750430740468"
$ws.Range("C1").Value = "0.9935404298917957"
$ws.Range("D1").Value = "0.9995915976099383"
$ws.Range("E1").Value = "0.9999996732049408"
$ws.Range("A2").Value = "0.5158054877847658"
$ws.Range("B2").Value = "0.9569204571328191"
$ws.Range("C2").Value = "0.9949509219699387"
$ws.Range("D2").Value = "0.9997494041500106"
$ws.Range("E2").Value = "0.9999998705981143"
$ws.Range("A3").Value = "0.5066973956625889"
$ws.Range("B3").Value = "0.9517112925674315"
$ws.Range("C3").Value = "0.9953900700852497"
$ws.Range("D3").Value = "0.9997887990788891"
$ws.Range("E3").Value = "0.9999998484275074"
$ws.Range("A4").Value = "0.5043441205810106"
$ws.Range("B4").Value = "0.949951866238407"
$ws.Range("C4").Value = "0.9951855792042685"
$ws.Range("D4").Value = "0.9998324788456466"
$ws.Range("E4").Value = "0.9999999292953776"
$ws.Range("A5").Value = "0.5046621137639407"
$ws.Range("B5").Value = "0.9509282401531843"
$ws.Range("C5").Value = "0.9946905471162312"
$ws.Range("D5").Value = "0.9998283016044427"
$ws.Range("E5").Value = "0.9999999443158842"
$ws.Range("A6").Value = "0.5028753796909075"
$ws.Range("B6").Value = "0.9521548411290554"
$ws.Range("C6").Value = "0.9949823254973134"
$ws.Range("D6").Value = "0.9998189019473729"
$ws.Range("E6").Value = "0.9999999469827321"
$ws.Range("A7").Value = "0.5017464050989884"
$ws.Range("B7").Value = "0.9518237955844949"
$ws.Range("C7").Value = "0.9948808300856948"
$ws.Range("D7").Value = "0.9998311666357742"
$ws.Range("E7").Value = "0.999999934526085"
$ws.Range("A8").Value = "0.5061910602324473"
$ws.Range("B8").Value = "0.9512385494882523"
$ws.Range("C8").Value = "0.9952109328298336"
$ws.Range("D8").Value = "0.9998052694530802"
$ws.Range("E8").Value = "0.9999999035841519"
$ws.Range("A9").Value = "0.5112106602408916"
$ws.Range("B9").Value = "0.9559861689750272"
$ws.Range("C9").Value = "0.9949761371894336"
$ws.Range("D9").Value = "0.9997627551198304"
$ws.Range("E9").Value = "0.9999998671778866"
$ws.Range("A10").Value = "0.545148932189429"
$ws.Range("B10").Value = "0.9311619612766192"
$ws.Range("C10").Value = "0.9938775975353399"
$ws.Range("D10").Value = "0.9995894712881097"
$ws.Range("E10").Value = "0.9999996066470124"
